$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1124.0377
$ws.Range("J17").Value = 1124.0377
$ws.Range("L17").Value = 3372.1131
$ws.Range("N17").Value = -3708.1131
$ws.Range("H40").Value = 17882768
$ws.Range("I40").Value = 17228.5
$ws.Range("J40").Value = 35748308
$ws.Range("K40").Value = 17228.5
$ws.Range("L40").Value = 35748308
$ws.Range("M40").Value = -17053.5
$ws.Range("N40").Value = -35748658
$ws.Range("H51").Value = 5203.75
$ws.Range("J51").Value = 4999.1665
$ws.Range("L51").Value = 4999.1665
$ws.Range("N51").Value = -5967.1665
$ws.Range("H98").Value = 1196.6666
$ws.Range("J98").Value = 2000
$ws.Range("L98").Value = 2000
$ws.Range("N98").Value = -4996
$ws.Range("H116").Value = 36949412
$ws.Range("I116").Value = 33800796
$ws.Range("K116").Value = 33800796
$ws.Range("M116").Value = -33797354
$ws.Range("H122").Value = 1196.6666
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 105846.01
$ws.Range("I132").Value = 166643.3
$ws.Range("J132").Value = 21665.154
$ws.Range("K132").Value = 499929.9
$ws.Range("L132").Value = 64995.462
$ws.Range("M132").Value = -497399.9
$ws.Range("N132").Value = -70055.462
$ws.Range("H138").Value = 4897.1694
$ws.Range("J138").Value = 5371.7017
$ws.Range("L138").Value = 16115.1051
$ws.Range("N138").Value = -26395.1051
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5570.48
$ws.Range("I32").Value = 3151.0256
$ws.Range("J32").Value = 14148.546
$ws.Range("K32").Value = 3151.0256
$ws.Range("L32").Value = 14148.546
$ws.Range("M32").Value = -2864.0256
$ws.Range("N32").Value = -14722.546
$ws.Range("H74").Value = 8622398
$ws.Range("I74").Value = 8930198
$ws.Range("K74").Value = 8930198
$ws.Range("M74").Value = -8929324
$ws.Range("H77").Value = 8622398
$ws.Range("I77").Value = 8930198
$ws.Range("K77").Value = 44650990
$ws.Range("M77").Value = -44646622
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H132").Value = 4370.732
$ws.Range("I132").Value = 1493.24
$ws.Range("K132").Value = 4479.72
$ws.Range("M132").Value = -1949.72
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5999.5
$ws.Range("I105").Value = 5999.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5999.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -4252.5
$ws.Range("N105").ClearContents()
$ws.Range("H119").Value = 55486.25
$ws.Range("J119").Value = 55486.25
$ws.Range("L119").Value = 55486.25
$ws.Range("N119").Value = -65162.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 69198.55499999999
$ws.Range("I19").Value = 4156.2
$ws.Range("J19").Value = 150501.5
$ws.Range("K19").Value = 4156.2
$ws.Range("L19").Value = 150501.5
$ws.Range("M19").Value = -3986.2
$ws.Range("N19").Value = -150841.5
$ws.Range("H24").Value = 69198.55499999999
$ws.Range("I24").Value = 4156.2
$ws.Range("J24").Value = 150501.5
$ws.Range("K24").Value = 4156.2
$ws.Range("L24").Value = 150501.5
$ws.Range("M24").Value = -3986.2
$ws.Range("N24").Value = -150841.5
$ws.Range("H31").Value = 4673
$ws.Range("I31").Value = 1511
$ws.Range("J31").Value = 5906.951
$ws.Range("K31").Value = 1511
$ws.Range("L31").Value = 5906.951
$ws.Range("M31").Value = -1216
$ws.Range("N31").Value = -6496.951
$ws.Range("H34").Value = 4673
$ws.Range("I34").Value = 1511
$ws.Range("J34").Value = 5906.951
$ws.Range("K34").Value = 1511
$ws.Range("L34").Value = 5906.951
$ws.Range("M34").Value = -1309
$ws.Range("N34").Value = -6310.951
$ws.Range("H62").Value = 46129.285
$ws.Range("I62").Value = 2952.5
$ws.Range("J62").Value = 63400
$ws.Range("K62").Value = 2952.5
$ws.Range("L62").Value = 63400
$ws.Range("M62").Value = -2328.5
$ws.Range("N62").Value = -64648
$ws.Range("H65").Value = 46129.285
$ws.Range("I65").Value = 2952.5
$ws.Range("J65").Value = 63400
$ws.Range("K65").Value = 14762.5
$ws.Range("L65").Value = 317000
$ws.Range("M65").Value = -11642.5
$ws.Range("N65").Value = -323240
$ws.Range("H134").Value = 2122.6785
$ws.Range("I134").Value = 2143.0417
$ws.Range("K134").Value = 6429.125100000001
$ws.Range("M134").Value = -3894.125100000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 2128.2856
$ws.Range("I33").Value = 1400
$ws.Range("J33").Value = 2419.6
$ws.Range("K33").Value = 1400
$ws.Range("L33").Value = 2419.6
$ws.Range("M33").Value = -1148
$ws.Range("N33").Value = -2923.6
$ws.Range("H102").Value = 31261950
$ws.Range("I102").Value = 55569788
$ws.Range("K102").Value = 55569788
$ws.Range("M102").Value = -55568166
$ws.Range("H113").Value = 1103.75
$ws.Range("I113").Value = 1184
$ws.Range("J113").Value = 1023.5
$ws.Range("K113").Value = 1184
$ws.Range("L113").Value = 1023.5
$ws.Range("M113").Value = 986
$ws.Range("N113").Value = -5363.5
$ws.Range("H122").Value = 428011.8
$ws.Range("I122").Value = 614267.25
$ws.Range("K122").Value = 1842801.75
$ws.Range("M122").Value = -1840351.75
$ws.Range("H126").Value = 4793.8
$ws.Range("I126").Value = 2799.375
$ws.Range("K126").Value = 8398.125
$ws.Range("M126").Value = -5928.125
$ws.Range("H132").Value = 94704.77
$ws.Range("I132").Value = 136918.8
$ws.Range("J132").Value = 4246.143
$ws.Range("K132").Value = 410756.4
$ws.Range("L132").Value = 12738.429
$ws.Range("M132").Value = -408226.4
$ws.Range("N132").Value = -17798.429
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5856.5713
$ws.Range("I7").Value = 1998.5
$ws.Range("K7").Value = 1998.5
$ws.Range("M7").Value = -1886.5
$ws.Range("H40").Value = 23813382
$ws.Range("I40").Value = 3571.7144
$ws.Range("K40").Value = 3571.7144
$ws.Range("M40").Value = -3435.7144
$ws.Range("H46").Value = 5999.2256
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("M46").Value = -4812
$ws.Range("H76").Value = 47666
$ws.Range("J76").Value = 47666
$ws.Range("L76").Value = 47666
$ws.Range("N76").Value = -48342
$ws.Range("H79").Value = 47666
$ws.Range("J79").Value = 47666
$ws.Range("L79").Value = 47666
$ws.Range("N79").Value = -50006
$ws.Range("H93").Value = 1707.7142
$ws.Range("I93").Value = 1295.8
$ws.Range("K93").Value = 1295.8
$ws.Range("M93").Value = -47.79999999999995
$ws.Range("H122").Value = 54427252
$ws.Range("I122").Value = 142861310
$ws.Range("K122").Value = 428583930
$ws.Range("M122").Value = -428581480
$ws.Range("H126").Value = 5856.5713
$ws.Range("I126").Value = 1998.5
$ws.Range("K126").Value = 5995.5
$ws.Range("M126").Value = -3525.5
$ws.Range("H132").Value = 3924.2188
$ws.Range("I132").Value = 3392.5186
$ws.Range("K132").Value = 10177.5558
$ws.Range("M132").Value = -7647.5558
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6859
$ws.Range("I2").Value = 7197
$ws.Range("K2").Value = 7197
$ws.Range("M2").Value = -7085
$ws.Range("H62").Value = 7859.222
$ws.Range("I62").Value = 3766.6
$ws.Range("J62").Value = 12975
$ws.Range("K62").Value = 3766.6
$ws.Range("L62").Value = 12975
$ws.Range("M62").Value = -3142.6
$ws.Range("N62").Value = -14223
$ws.Range("H65").Value = 7859.222
$ws.Range("I65").Value = 3766.6
$ws.Range("J65").Value = 12975
$ws.Range("K65").Value = 18833
$ws.Range("L65").Value = 64875
$ws.Range("M65").Value = -15713
$ws.Range("N65").Value = -71115
$ws.Range("H126").Value = 3899.8333
$ws.Range("I126").Value = 4250
$ws.Range("K126").Value = 12750
$ws.Range("M126").Value = -10280
$ws.Range("H132").Value = 24160332
$ws.Range("I132").Value = 4630414.5
$ws.Range("J132").Value = 45465696
$ws.Range("K132").Value = 13891243.5
$ws.Range("L132").Value = 136397088
$ws.Range("M132").Value = -13888713.5
$ws.Range("N132").Value = -136402148
$ws.Range("H136").Value = 8549.117
$ws.Range("I136").Value = 2149.8572
$ws.Range("J136").Value = 11408.361
$ws.Range("K136").Value = 6449.571599999999
$ws.Range("L136").Value = 34225.083
$ws.Range("M136").Value = -3899.571599999999
$ws.Range("N136").Value = -39325.083
